$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 563.25
$ws.Range("I2").Value = 459.88235
$ws.Range("J2").Value = 814.2857
$ws.Range("K2").Value = 459.88235
$ws.Range("L2").Value = 814.2857
$ws.Range("M2").Value = -346.88235
$ws.Range("N2").Value = -1040.2857
$ws.Range("H45").Value = 1280
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1280
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1280
$ws.Range("N45").Value = -2034
$ws.Range("M45").ClearContents()
$ws.Range("H61").Value = 1768.9375
$ws.Range("I61").Value = 1527.7931
$ws.Range("K61").Value = 1527.7931
$ws.Range("M61").Value = -1315.7931
$ws.Range("H116").Value = 563.25
$ws.Range("I116").Value = 459.88235
$ws.Range("J116").Value = 814.2857
$ws.Range("K116").Value = 459.88235
$ws.Range("L116").Value = 814.2857
$ws.Range("M116").Value = 1834.11765
$ws.Range("N116").Value = -5402.2857
$ws.Range("H132").Value = 2593.6736
$ws.Range("I132").Value = 1990.0488
$ws.Range("J132").Value = 5687.25
$ws.Range("K132").Value = 5970.1464
$ws.Range("L132").Value = 17061.75
$ws.Range("M132").Value = -3440.1464
$ws.Range("N132").Value = -22121.75
$ws.Range("H135").Value = 40936
$ws.Range("J135").Value = 40936
$ws.Range("L135").Value = 40936
$ws.Range("N135").Value = -51076
$ws.Range("H136").Value = 1768.9375
$ws.Range("I136").Value = 1527.7931
$ws.Range("K136").Value = 4583.379300000001
$ws.Range("M136").Value = -2033.379300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 563.25
$ws.Range("I3").Value = 459.88235
$ws.Range("J3").Value = 814.2857
$ws.Range("K3").Value = 459.88235
$ws.Range("L3").Value = 814.2857
$ws.Range("M3").Value = -345.88235
$ws.Range("N3").Value = -1042.2857
$ws.Range("H134").Value = 2511.4084
$ws.Range("I134").Value = 1534.2424
$ws.Range("J134").Value = 3360
$ws.Range("K134").Value = 4602.7272
$ws.Range("L134").Value = 10080
$ws.Range("M134").Value = -2067.7272
$ws.Range("N134").Value = -15150
$ws.Range("H135").Value = 37544.445
$ws.Range("J135").Value = 37544.445
$ws.Range("L135").Value = 37544.445
$ws.Range("N135").Value = -47684.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 27780534
$ws.Range("I16").Value = 111111110
$ws.Range("J16").Value = 3675.3333
$ws.Range("K16").Value = 111111110
$ws.Range("L16").Value = 3675.3333
$ws.Range("M16").Value = -111110823
$ws.Range("N16").Value = -4249.3333
$ws.Range("H86").Value = 2856.3635
$ws.Range("I86").Value = 2814.25
$ws.Range("K86").Value = 2814.25
$ws.Range("M86").Value = -1691.25
$ws.Range("H89").Value = 2856.3635
$ws.Range("I89").Value = 2814.25
$ws.Range("K89").Value = 14071.25
$ws.Range("M89").Value = -8455.25
$ws.Range("H113").Value = 27780534
$ws.Range("I113").Value = 111111110
$ws.Range("J113").Value = 3675.3333
$ws.Range("K113").Value = 111111110
$ws.Range("L113").Value = 3675.3333
$ws.Range("M113").Value = -111108940
$ws.Range("N113").Value = -8015.3333
$ws.Range("H122").Value = 1988.8695
$ws.Range("I122").Value = 1327.6875
$ws.Range("K122").Value = 3983.0625
$ws.Range("M122").Value = -1533.0625
$ws.Range("H141").Value = 25700
$ws.Range("J141").Value = 25700
$ws.Range("L141").Value = 25700
$ws.Range("N141").Value = -36060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 766.2273
$ws.Range("I113").Value = 665.9167
$ws.Range("J113").Value = 886.6
$ws.Range("K113").Value = 1997.7501
$ws.Range("L113").Value = 2659.8
$ws.Range("M113").Value = 172.2499
$ws.Range("N113").Value = -6999.8
$ws.Range("H132").Value = 1808.7368
$ws.Range("I132").Value = 792.8570999999999
$ws.Range("J132").Value = 2401.3333
$ws.Range("K132").Value = 7135.7139
$ws.Range("L132").Value = 21611.9997
$ws.Range("M132").Value = -4605.7139
$ws.Range("N132").Value = -26671.9997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5714.0713
$ws.Range("I7").Value = 3000.75
$ws.Range("J7").Value = 6799.4
$ws.Range("K7").Value = 3000.75
$ws.Range("L7").Value = 6799.4
$ws.Range("M7").Value = -2888.75
$ws.Range("N7").Value = -7023.4
$ws.Range("H122").Value = 6368.579
$ws.Range("I122").Value = 4000.3
$ws.Range("K122").Value = 12000.9
$ws.Range("M122").Value = -9550.900000000001
$ws.Range("H126").Value = 5714.0713
$ws.Range("I126").Value = 3000.75
$ws.Range("J126").Value = 6799.4
$ws.Range("K126").Value = 9002.25
$ws.Range("L126").Value = 20398.2
$ws.Range("M126").Value = -6532.25
$ws.Range("N126").Value = -25338.2
$ws.Range("H132").Value = 3019.0967
$ws.Range("I132").Value = 1653.7028
$ws.Range("J132").Value = 5039.88
$ws.Range("K132").Value = 4961.1084
$ws.Range("L132").Value = 15119.64
$ws.Range("M132").Value = -2431.1084
$ws.Range("N132").Value = -20179.64

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H113").Value = 358
$ws.Range("I113").Value = 312.2143
$ws.Range("K113").Value = 936.6428999999999
$ws.Range("M113").Value = 1233.3571
$ws.Range("H115").Value = 26021.738
$ws.Range("J115").Value = 26021.738
$ws.Range("L115").Value = 26021.738
$ws.Range("N115").Value = -29155.738
$ws.Range("H122").Value = 3215.8064
$ws.Range("I122").Value = 2074.25
$ws.Range("J122").Value = 5291.364
$ws.Range("K122").Value = 6222.75
$ws.Range("L122").Value = 15874.092
$ws.Range("M122").Value = -3772.75
$ws.Range("N122").Value = -20774.092
$ws.Range("H126").Value = 2260.9119
$ws.Range("I126").Value = 1770.7646
$ws.Range("J126").Value = 2751.0588
$ws.Range("K126").Value = 5312.293799999999
$ws.Range("L126").Value = 8253.1764
$ws.Range("M126").Value = -2842.293799999999
$ws.Range("N126").Value = -13193.1764
$ws.Range("H132").Value = 10102860
$ws.Range("I132").Value = 896.7143
$ws.Range("J132").Value = 27781296
$ws.Range("K132").Value = 2690.1429
$ws.Range("L132").Value = 83343888
$ws.Range("M132").Value = -160.1428999999998
$ws.Range("N132").Value = -83348948

